# Insert a new "spatial_entropy" column before the existing "diversity_shannon"
# column (currently column D), shifting subsequent columns one to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; this shifts D:G -> E:H.
$ws.Range("D1").EntireColumn.Insert()

# Header for the newly inserted column - match the bold/centered/bordered
# header style used by the rest of row 1.
$ws.Range("D1").Value = "spatial_entropy"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("D1").VerticalAlignment = -4160     # xlTop
$ws.Range("D1").Borders.LineStyle = 1

# New spatial_entropy values for each data row
$ws.Range("D2").Value = 3.681212142935576
$ws.Range("D3").Value = 3.720209773733151
$ws.Range("D4").Value = 3.508339164031677
$ws.Range("D5").Value = 3.424244434780839

# Restore the original selection (A1) now that the sheet spans A1:H5.
$ws.Range("A1").Select()
